$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Existing MRF TC: refresh the "Date" column (B2:B48) with the
# timestamps produced by the latest Katalon test run (Fri Oct 25 2024).
$dates = @(
    "Fri Oct 25 11:35:29 EDT 2024",
    "Fri Oct 25 11:35:39 EDT 2024",
    "Fri Oct 25 11:35:49 EDT 2024",
    "Fri Oct 25 11:36:00 EDT 2024",
    "Fri Oct 25 11:36:10 EDT 2024",
    "Fri Oct 25 11:36:20 EDT 2024",
    "Fri Oct 25 11:36:30 EDT 2024",
    "Fri Oct 25 11:36:40 EDT 2024",
    "Fri Oct 25 11:36:50 EDT 2024",
    "Fri Oct 25 11:37:00 EDT 2024",
    "Fri Oct 25 11:37:09 EDT 2024",
    "Fri Oct 25 11:37:19 EDT 2024",
    "Fri Oct 25 11:37:29 EDT 2024",
    "Fri Oct 25 11:37:39 EDT 2024",
    "Fri Oct 25 11:37:49 EDT 2024",
    "Fri Oct 25 11:37:59 EDT 2024",
    "Fri Oct 25 11:38:09 EDT 2024",
    "Fri Oct 25 11:38:19 EDT 2024",
    "Fri Oct 25 11:38:29 EDT 2024",
    "Fri Oct 25 11:38:39 EDT 2024",
    "Fri Oct 25 11:38:49 EDT 2024",
    "Fri Oct 25 11:38:59 EDT 2024",
    "Fri Oct 25 11:39:08 EDT 2024",
    "Fri Oct 25 11:39:18 EDT 2024",
    "Fri Oct 25 11:39:28 EDT 2024",
    "Fri Oct 25 11:39:37 EDT 2024",
    "Fri Oct 25 11:39:47 EDT 2024",
    "Fri Oct 25 11:39:57 EDT 2024",
    "Fri Oct 25 11:40:07 EDT 2024",
    "Fri Oct 25 11:40:17 EDT 2024",
    "Fri Oct 25 11:40:27 EDT 2024",
    "Fri Oct 25 11:40:37 EDT 2024",
    "Fri Oct 25 11:40:47 EDT 2024",
    "Fri Oct 25 11:40:57 EDT 2024",
    "Fri Oct 25 11:41:07 EDT 2024",
    "Fri Oct 25 11:41:17 EDT 2024",
    "Fri Oct 25 11:41:27 EDT 2024",
    "Fri Oct 25 11:41:37 EDT 2024",
    "Fri Oct 25 11:41:47 EDT 2024",
    "Fri Oct 25 11:41:57 EDT 2024",
    "Fri Oct 25 11:42:07 EDT 2024",
    "Fri Oct 25 11:42:17 EDT 2024",
    "Fri Oct 25 11:42:26 EDT 2024",
    "Fri Oct 25 11:42:36 EDT 2024",
    "Fri Oct 25 11:42:46 EDT 2024",
    "Fri Oct 25 11:42:56 EDT 2024",
    "Fri Oct 25 11:43:06 EDT 2024"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}

